# Test data added for Italy
#
# Mirrors: right-click the "Slovakia" sheet tab -> Move or Copy -> Create a
# copy -> place after "Slovakia" -> rename the copy to "Italy" -> update the
# two market-specific cells (B4 "user story" id and B2 "description/market"
# label) on the new sheet.

$wb = $excel.ActiveWorkbook

# Duplicate "Slovakia" and drop the copy right after it (becomes the last tab).
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $slovakia)

# The copy is now the last sheet and is the active sheet, just like Excel
# leaves it after "Move or Copy -> Create a copy".
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Update the copied sheet's content for the new market. Write B4 before B2 so
# new shared-string entries are appended in the same order Excel would use
# (B4's string first, then B2's).
$italy.Range("B4").Value = "NGC-3145/T2159 "
$italy.Range("B2").Value = "Italy Market"
